# Updated symbol list on Thu Jan 12 05:51:43 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns on the crypto listing sheet.
# Values are stored as text (they include things like "%" and keep exact
# trailing-zero formatting), so we force the target cells to Text format
# before writing the new value - this mirrors how the original sheet keeps
# these columns as plain text rather than numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2";  Value = "285.54" },
    @{ Cell = "E2";  Value = "3.46%" },
    @{ Cell = "D3";  Value = "28.87" },
    @{ Cell = "E3";  Value = "6.40%" },
    @{ Cell = "D4";  Value = "4.932" },
    @{ Cell = "E4";  Value = "1.54%" },
    @{ Cell = "D5";  Value = "0.06507" },
    @{ Cell = "E5";  Value = "1.68%" },
    @{ Cell = "D6";  Value = "7.233" },
    @{ Cell = "E6";  Value = "4.17%" },
    @{ Cell = "D7";  Value = "1.348" },
    @{ Cell = "E7";  Value = "12.54%" },
    @{ Cell = "D8";  Value = "0.9109" },
    @{ Cell = "E8";  Value = "3.72%" },
    @{ Cell = "D9";  Value = "0.1552" },
    @{ Cell = "E9";  Value = "2.19%" },
    @{ Cell = "D10"; Value = "0.06589" },
    @{ Cell = "E10"; Value = "29.25%" },
    @{ Cell = "D11"; Value = "0.07723" },
    @{ Cell = "E11"; Value = "2.76%" },
    @{ Cell = "D12"; Value = "0.02989" },
    @{ Cell = "E12"; Value = "0.48%" },
    @{ Cell = "E13"; Value = "-0.12%" },
    @{ Cell = "D14"; Value = "0.001604" },
    @{ Cell = "E14"; Value = "2.70%" },
    @{ Cell = "D15"; Value = "0.0006534" },
    @{ Cell = "E15"; Value = "2.30%" },
    @{ Cell = "D16"; Value = "0.006028" },
    @{ Cell = "E16"; Value = "-2.53%" },
    @{ Cell = "D18"; Value = "3.390" },
    @{ Cell = "E18"; Value = "2.52%" },
    @{ Cell = "D19"; Value = "2.240" },
    @{ Cell = "E19"; Value = "-0.59%" },
    @{ Cell = "D20"; Value = "0.3148" },
    @{ Cell = "E20"; Value = "0.40%" },
    @{ Cell = "E21"; Value = "0.14%" },
    @{ Cell = "D22"; Value = "4.042" },
    @{ Cell = "E22"; Value = "3.01%" },
    @{ Cell = "D23"; Value = "0.1556" },
    @{ Cell = "E23"; Value = "12.74%" },
    @{ Cell = "D24"; Value = "0.04491" },
    @{ Cell = "E24"; Value = "1.67%" },
    @{ Cell = "D25"; Value = "0.001192" },
    @{ Cell = "E25"; Value = "1.36%" },
    @{ Cell = "D26"; Value = "0.004325" },
    @{ Cell = "E26"; Value = "11.97%" },
    @{ Cell = "D28"; Value = "0.0001184" },
    @{ Cell = "E28"; Value = "-1.31%" },
    @{ Cell = "D29"; Value = "0.0001637" },
    @{ Cell = "E29"; Value = "-15.63%" },
    @{ Cell = "D40"; Value = "0.04158" },
    @{ Cell = "E40"; Value = "-0.07%" },
    @{ Cell = "D41"; Value = "0.006717" },
    @{ Cell = "E41"; Value = "-1.38%" },
    @{ Cell = "E42"; Value = "5.18%" },
    @{ Cell = "D43"; Value = "0.002188" },
    @{ Cell = "E43"; Value = "1.77%" },
    @{ Cell = "D44"; Value = "0.01175" },
    @{ Cell = "E44"; Value = "-0.94%" },
    @{ Cell = "D45"; Value = "0.00005445" },
    @{ Cell = "E45"; Value = "3.21%" },
    @{ Cell = "D46"; Value = "1.562" },
    @{ Cell = "E46"; Value = "-7.03%" },
    @{ Cell = "D47"; Value = "0.01852" },
    @{ Cell = "E47"; Value = "0.15%" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    # Force text storage so numeric-looking strings (e.g. "285.54") and
    # percent-looking strings (e.g. "3.46%") are not reinterpreted by Excel
    # as numbers/percentages, matching the original column formatting.
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
